# Apply updated simulation results across the workbook.
# Sheet order: 1 = Coefficients, 2 = Performance Metrics,
#              3 = Runtime Analysis, 4 = Memory Usage

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Coefficients - update "Encrypted Coefficient" column (C2:C22)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C2").Value  = -17854.32704906524
$ws1.Range("C3").Value  = 14042.59545713462
$ws1.Range("C4").Value  = 5513.877038814026
$ws1.Range("C5").Value  = 4701.117548934597
$ws1.Range("C6").Value  = 4363.940245938534
$ws1.Range("C7").Value  = -4350.72552279905
$ws1.Range("C8").Value  = -2201.599338929234
$ws1.Range("C9").Value  = -1364.822336431141
$ws1.Range("C10").Value = -427.661478598414
$ws1.Range("C11").Value = -321.2333687677409
$ws1.Range("C12").Value = -201.7694550906396
$ws1.Range("C13").Value = 201.2629543625523
$ws1.Range("C14").Value = 185.0718710626106
$ws1.Range("C15").Value = 151.5601082696767
$ws1.Range("C16").Value = 143.0183857675093
$ws1.Range("C17").Value = 114.7607661645206
$ws1.Range("C18").Value = 10.59737628136645
$ws1.Range("C19").Value = -7.089485865492861
$ws1.Range("C20").Value = 6.209684236669091
$ws1.Range("C21").Value = -0.2926168162575777
$ws1.Range("C22").Value = -0.000000001064549337570497

# ---------------------------------------------------------------------------
# Sheet 2: Performance Metrics - update "Encrypted Inference" column (C2:C5)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C2").Value = 2.473501344543386
$ws2.Range("C3").Value = 15.0277061464848
$ws2.Range("C4").Value = 3.876558544183848
$ws2.Range("C5").Value = 0.9263891185317032

# ---------------------------------------------------------------------------
# Sheet 3: Runtime Analysis - insert "Encryption Time" / "Decryption Time"
# rows, rename "Total without preprocessing" and refresh every value.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Insert two new blank rows right after "Training Time" (row 4), pushing the
# old Evaluation Time / Total without preprocessing / Total Runtime rows down.
$ws3.Rows.Item(5).Insert()
$ws3.Rows.Item(5).Insert()

# Row 2: Context Setup Time
$ws3.Range("C2").Value = 0.4714062213897705

# Row 3: Preprocessing Time
$ws3.Range("B3").Value = 3.666933059692383
$ws3.Range("C3").Value = 5.661703824996948

# Row 4: Training Time
$ws3.Range("B4").Value = 0.004482507705688477
$ws3.Range("C4").Value = 0.004482507705688477

# Row 5 (new): Encryption Time
$ws3.Range("A5").Value = "Encryption Time"
$ws3.Range("C5").Value = 4.491698265075684

# Row 6 (new): Decryption Time
$ws3.Range("A6").Value = "Decryption Time"
$ws3.Range("C6").Value = 0.6985993385314941

# Row 7 (was row 5): Evaluation Time
$ws3.Range("A7").Value = "Evaluation Time"
$ws3.Range("B7").Value = 0
$ws3.Range("C7").Value = 14.0158166885376

# Row 8 (was row 6): Total without initial preprocessing (renamed)
$ws3.Range("A8").Value = "Total without initial preprocessing"
$ws3.Range("B8").Value = 0.004482507705688477
$ws3.Range("C8").Value = 19.67752051353455

# Row 9 (was row 7): Total Runtime
$ws3.Range("A9").Value = "Total Runtime"
$ws3.Range("B9").Value = 3.671415567398071
$ws3.Range("C9").Value = 23.34893608093262

# ---------------------------------------------------------------------------
# Sheet 4: Memory Usage - update "Encrypted Inference (KB)" column (C2:C4)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("C2").Value = 326.1572265625
$ws4.Range("C3").Value = 229.923828125
$ws4.Range("C4").Value = 258581.541015625
